$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D2", "D3", "D4", "D5", "D6", "D9", "D12", "D13", "D14", "D15", "D17", "D18", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D34", "D36", "D37", "D38", "D39", "D42", "D43", "D45", "D47", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '67.879.09'
$ws.Range('E2').Value = '  +1.15%  '
$ws.Range('D3').Value = '2.508.91'
$ws.Range('E3').Value = '  +0.93%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = '589.34'
$ws.Range('E5').Value = '  +1.00%  '
$ws.Range('D6').Value = '177.16'
$ws.Range('E6').Value = '  +3.42%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +0.61%  '
$ws.Range('D9').Value = '0.141'
$ws.Range('E9').Value = '  +3.41%  '
$ws.Range('E10').Value = '  -0.27%  '
$ws.Range('E11').Value = '  +2.55%  '
$ws.Range('D12').Value = '4.96'
$ws.Range('E12').Value = '  +0.69%  '
$ws.Range('D13').Value = '25.74'
$ws.Range('E13').Value = '  +1.56%  '
$ws.Range('D14').Value = '2.926.23'
$ws.Range('D15').Value = '67.841.65'
$ws.Range('E15').Value = '  +1.16%  '
$ws.Range('E16').Value = '  +1.34%  '
$ws.Range('D17').Value = '2.495.21'
$ws.Range('E17').Value = '  +1.33%  '
$ws.Range('D18').Value = '11.02'
$ws.Range('E18').Value = '  +0.23%  '
$ws.Range('E19').Value = '  +2.24%  '
$ws.Range('D20').Value = '353.28'
$ws.Range('E20').Value = '  +1.53%  '
$ws.Range('D21').Value = '4.05'
$ws.Range('E21').Value = '  +0.61%  '
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').Value = '70.83'
$ws.Range('E23').Value = '  +3.41%  '
$ws.Range('D24').Value = '4.30'
$ws.Range('E24').Value = '  +1.93%  '
$ws.Range('D25').Value = '1.76'
$ws.Range('E25').Value = '  -1.05%  '
$ws.Range('D26').Value = '9.15'
$ws.Range('E26').Value = '  -1.21%  '
$ws.Range('D27').Value = '2.593.98'
$ws.Range('E27').Value = '  -0.81%  '
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').Value = '  -0.12%  '
$ws.Range('D29').Value = '0.0₃0919'
$ws.Range('E29').Value = '  +1.94%  '
$ws.Range('D30').Value = '508.06'
$ws.Range('E30').Value = '  -0.22%  '
$ws.Range('D31').Value = '7.87'
$ws.Range('E31').Value = '  +1.48%  '
$ws.Range('E32').Value = '  +2.90%  '
$ws.Range('E33').Value = '  +1.01%  '
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('E35').Value = '  +5.13%  '
$ws.Range('D36').Value = '164.57'
$ws.Range('E36').Value = '  +2.96%  '
$ws.Range('D37').Value = '18.42'
$ws.Range('E37').Value = '  +1.13%  '
$ws.Range('D38').Value = '18.65'
$ws.Range('E38').Value = '  -0.20%  '
$ws.Range('D39').Value = '1.35'
$ws.Range('E39').Value = '  +0.94%  '
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('E41').Value = '  +3.18%  '
$ws.Range('D42').Value = '4.90'
$ws.Range('E42').Value = '  +1.90%  '
$ws.Range('D43').Value = '0.330'
$ws.Range('E43').Value = '  +0.67%  '
$ws.Range('E44').Value = '  +5.55%  '
$ws.Range('D45').Value = '145.80'
$ws.Range('E45').Value = '  +2.45%  '
$ws.Range('E46').Value = '  +3.12%  '
$ws.Range('D47').Value = '0.520'
$ws.Range('E47').Value = '  +1.28%  '
$ws.Range('E48').Value = '  +3.76%  '
$ws.Range('D49').Value = '0.0744'
$ws.Range('E49').Value = '  +1.72%  '
$ws.Range('E50').Value = '  +1.91%  '
$ws.Range('D51').Value = '0.588'
$ws.Range('E51').Value = '  +1.02%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
